# Insert a new data row at row 211 (Betarraga / Feria Lagunitas de Puerto Montt
# weekly record), pushing the existing rows 211-294 down to 212-295.
# The new row duplicates the row that used to be at 211 (same volume/price/
# unit/origin data) but is stamped with a new "Fecha" (date) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 211; everything from old row 211 onward shifts to 212+.
$ws.Rows("211:211").Insert()

# Populate the newly blank row 211 with the same data that is now sitting in
# row 212 (the row that used to be 211, shifted down by the insert above).
$ws.Range("A211:R211").Value = $ws.Range("A212:R212").Value()

# Stamp the new row with its own date (Fecha column D) -> serial 44704
# (2022-05-23), matching the author's edit.
$ws.Range("D211").Value = "5/23/2022"
